$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D:F) for City / State / Zipcode. This shifts the
# existing Phone / Phone Tag / Notes / Relationships columns from D:G to G:J.
$ws.Columns("D:F").Insert()

# Notes column (now I) - first newly-introduced value so it becomes the
# first new shared string.
$ws.Cells.Item(2, 9).Value = "Huge Organization"

# New column headers.
$ws.Cells.Item(1, 4).Value = "City"
$ws.Cells.Item(1, 5).Value = "State"
$ws.Cells.Item(1, 6).Value = "Zipcode"

# Austin, TX rows (2-8).
$ws.Cells.Item(2, 4).Value = "Austin"
$ws.Cells.Item(2, 5).Value = "TX"
$ws.Cells.Item(2, 6).Value = 78701

$ws.Cells.Item(3, 4).Value = "Austin"
$ws.Cells.Item(3, 5).Value = "TX"
$ws.Cells.Item(3, 6).Value = 78702

$ws.Cells.Item(4, 4).Value = "Austin"
$ws.Cells.Item(4, 5).Value = "TX"
$ws.Cells.Item(4, 6).Value = 78703

$ws.Cells.Item(5, 4).Value = "Austin"
$ws.Cells.Item(5, 5).Value = "TX"
$ws.Cells.Item(5, 6).Value = 78704

$ws.Cells.Item(6, 4).Value = "Austin"
$ws.Cells.Item(6, 5).Value = "TX"
$ws.Cells.Item(6, 6).Value = 78705

$ws.Cells.Item(7, 4).Value = "Austin"
$ws.Cells.Item(7, 5).Value = "TX"
$ws.Cells.Item(7, 6).Value = 78706

$ws.Cells.Item(8, 4).Value = "Austin"
$ws.Cells.Item(8, 5).Value = "TX"
$ws.Cells.Item(8, 6).Value = 78707

# New York, NY rows (9-12).
$ws.Cells.Item(9, 4).Value = "New York"
$ws.Cells.Item(9, 5).Value = "NY"
$ws.Cells.Item(9, 6).Value = 800001

$ws.Cells.Item(10, 4).Value = "New York"
$ws.Cells.Item(10, 5).Value = "NY"
$ws.Cells.Item(10, 6).Value = 800002

$ws.Cells.Item(11, 4).Value = "New York"
$ws.Cells.Item(11, 5).Value = "NY"
$ws.Cells.Item(11, 6).Value = 800003

$ws.Cells.Item(12, 4).Value = "New York"
$ws.Cells.Item(12, 5).Value = "NY"
$ws.Cells.Item(12, 6).Value = 800004

# Column widths for the new City/State/Zipcode block and the Notes column.
$ws.Range("D1:F1").ColumnWidth = 17
$ws.Range("I1").ColumnWidth = 10.83

# Mailing-address display / print formatting improvements.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the sheet with row 2 selected, as in the authored workbook.
$ws.Rows("2:2").Select()
